$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 20,14
$arr[0,0] = 0.026992
$arr[0,1] = 0.08097599999999999
$arr[0,2] = 0.004182906599909731
$arr[0,3] = 0.00420788870005516
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.306376666666667
$arr[0,7] = 3.91913
$arr[0,8] = 0.06159635513812315
$arr[0,9] = 0.07271399171915481
$arr[0,10] = 0.03526171898666666
$arr[0,11] = 0.3173554708799999
$arr[0,12] = 0.000257651800437639
$arr[0,13] = 0.000305972384090936
$arr[1,0] = 0.026992
$arr[1,1] = 0.08097599999999999
$arr[1,2] = 0.004182906599909731
$arr[1,3] = 0.00420788870005516
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 8.648731
$arr[1,7] = 25.946193
$arr[1,8] = 0.4077922698431246
$arr[1,9] = 0.4813954277979023
$arr[1,10] = 0.233446547152
$arr[1,11] = 2.101018924368
$arr[1,12] = 0.001705756976918976
$arr[1,13] = 0.002025658380889013
$arr[2,0] = 0.026992
$arr[2,1] = 0.08097599999999999
$arr[2,2] = 0.004182906599909731
$arr[2,3] = 0.00420788870005516
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.637617
$arr[2,7] = 1.912851
$arr[2,8] = 0.03006398091472189
$arr[2,9] = 0.03549028273468269
$arr[2,10] = 0.017210558064
$arr[2,11] = 0.154895022576
$arr[2,12] = 0.0001257548241877504
$arr[2,13] = 0.0001493391596810341
$arr[3,0] = 0.026992
$arr[3,1] = 0.08097599999999999
$arr[3,2] = 0.004182906599909731
$arr[3,3] = 0.00420788870005516
$arr[3,4] = 2
$arr[3,5] = 1
$arr[3,6] = 9.728125
$arr[3,7] = 19.45625
$arr[3,8] = 0.4586862714388558
$arr[3,9] = 0.3609835859963323
$arr[3,10] = 0.26258155
$arr[3,11] = 1.5754893
$arr[3,12] = 0.001918641832089576
$arr[3,13] = 0.001518978752419357
$arr[4,0] = 0.026992
$arr[4,1] = 0.08097599999999999
$arr[4,2] = 0.004182906599909731
$arr[4,3] = 0.00420788870005516
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 0.8878186666666666
$arr[4,7] = 2.663456
$arr[4,8] = 0.0418611226651744
$arr[4,9] = 0.0494167117519279
$arr[4,10] = 0.02396400145066667
$arr[4,11] = 0.215676013056
$arr[4,12] = 0.0001751011662757888
$arr[4,13] = 0.0002079400229748205
$arr[5,0] = 6.300519666666666
$arr[5,1] = 18.901559
$arr[5,2] = 0.976381346197431
$arr[5,3] = 0.9822127115383066
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 1.306376666666667
$arr[5,7] = 3.91913
$arr[5,8] = 0.06159635513812315
$arr[5,9] = 0.07271399171915481
$arr[5,10] = 8.230851880407776
$arr[5,11] = 74.07766692366998
$arr[5,12] = 0.06014153215061573
$arr[5,13] = 0.07142060697324502
$arr[6,0] = 6.300519666666666
$arr[6,1] = 18.901559
$arr[6,2] = 0.976381346197431
$arr[6,3] = 0.9822127115383066
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 8.648731
$arr[6,7] = 25.946193
$arr[6,8] = 0.4077922698431246
$arr[6,9] = 0.4813954277979023
$arr[6,10] = 54.49149975720966
$arr[6,11] = 490.423497814887
$arr[6,12] = 0.3981607653983361
$arr[6,13] = 0.4728327084595207
$arr[7,0] = 6.300519666666666
$arr[7,1] = 18.901559
$arr[7,2] = 0.976381346197431
$arr[7,3] = 0.9822127115383066
$arr[7,4] = 2
$arr[7,5] = 0.6666666666666666
$arr[7,6] = 0.637617
$arr[7,7] = 1.912851
$arr[7,8] = 0.03006398091472189
$arr[7,9] = 0.03549028273468269
$arr[7,10] = 4.017318448300999
$arr[7,11] = 36.15586603470899
$arr[7,12] = 0.02935391015757003
$arr[7,13] = 0.03485900683809384
$arr[8,0] = 6.300519666666666
$arr[8,1] = 18.901559
$arr[8,2] = 0.976381346197431
$arr[8,3] = 0.9822127115383066
$arr[8,4] = 2
$arr[8,5] = 1
$arr[8,6] = 9.728125
$arr[8,7] = 19.45625
$arr[8,8] = 0.4586862714388558
$arr[8,9] = 0.3609835859963323
$arr[8,10] = 61.29224288229167
$arr[8,11] = 367.75345729375
$arr[8,12] = 0.4478527191897503
$arr[8,13] = 0.354562666822279
$arr[9,0] = 6.300519666666666
$arr[9,1] = 18.901559
$arr[9,2] = 0.976381346197431
$arr[9,3] = 0.9822127115383066
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 0.8878186666666666
$arr[9,7] = 2.663456
$arr[9,8] = 0.0418611226651744
$arr[9,9] = 0.0494167117519279
$arr[9,10] = 5.593718969767111
$arr[9,11] = 50.343470727904
$arr[9,12] = 0.04087241930115877
$arr[9,13] = 0.04853772244516801
$arr[10,0] = 0.1149325
$arr[10,1] = 0.229865
$arr[10,2] = 0.01781090370458377
$arr[10,3] = 0.01194485200600399
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 1.306376666666667
$arr[10,7] = 3.91913
$arr[10,8] = 0.06159635513812315
$arr[10,9] = 0.07271399171915481
$arr[10,10] = 0.1501451362416666
$arr[10,11] = 0.9008708174499999
$arr[10,12] = 0.001097086749918455
$arr[10,13] = 0.0008685578698511041
$arr[11,0] = 0.1149325
$arr[11,1] = 0.229865
$arr[11,2] = 0.01781090370458377
$arr[11,3] = 0.01194485200600399
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 8.648731
$arr[11,7] = 25.946193
$arr[11,8] = 0.4077922698431246
$arr[11,9] = 0.4813954277979023
$arr[11,10] = 0.9940202756575
$arr[11,11] = 5.964121653945
$arr[11,12] = 0.007263148849649532
$arr[11,13] = 0.005750197141412923
$arr[12,0] = 0.1149325
$arr[12,1] = 0.229865
$arr[12,2] = 0.01781090370458377
$arr[12,3] = 0.01194485200600399
$arr[12,4] = 2
$arr[12,5] = 0.6666666666666666
$arr[12,6] = 0.637617
$arr[12,7] = 1.912851
$arr[12,8] = 0.03006398091472189
$arr[12,9] = 0.03549028273468269
$arr[12,10] = 0.07328291585249999
$arr[12,11] = 0.4396974951149999
$arr[12,12] = 0.0005354666690485559
$arr[12,13] = 0.0004239261749170234
$arr[13,0] = 0.1149325
$arr[13,1] = 0.229865
$arr[13,2] = 0.01781090370458377
$arr[13,3] = 0.01194485200600399
$arr[13,4] = 2
$arr[13,5] = 1
$arr[13,6] = 9.728125
$arr[13,7] = 19.45625
$arr[13,8] = 0.4586862714388558
$arr[13,9] = 0.3609835859963323
$arr[13,10] = 1.1180777265625
$arr[13,11] = 4.47231090625
$arr[13,12] = 0.008169617011212034
$arr[13,13] = 0.004311895511322806
$arr[14,0] = 0.1149325
$arr[14,1] = 0.229865
$arr[14,2] = 0.01781090370458377
$arr[14,3] = 0.01194485200600399
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 0.8878186666666666
$arr[14,7] = 2.663456
$arr[14,8] = 0.0418611226651744
$arr[14,9] = 0.0494167117519279
$arr[14,10] = 0.1020392189066667
$arr[14,11] = 0.6122353134399999
$arr[14,12] = 0.0007455844247551903
$arr[14,13] = 0.0005902753085001371
$arr[15,0] = 0.010485
$arr[15,1] = 0.031455
$arr[15,2] = 0.001624843498075486
$arr[15,3] = 0.001634547755634201
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 1.306376666666667
$arr[15,7] = 3.91913
$arr[15,8] = 0.06159635513812315
$arr[15,9] = 0.07271399171915481
$arr[15,10] = 0.01369735935
$arr[15,11] = 0.12327623415
$arr[15,12] = 0.000100084437151328
$arr[15,13] = 0.0001188544919677484
$arr[16,0] = 0.010485
$arr[16,1] = 0.031455
$arr[16,2] = 0.001624843498075486
$arr[16,3] = 0.001634547755634201
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 8.648731
$arr[16,7] = 25.946193
$arr[16,8] = 0.4077922698431246
$arr[16,9] = 0.4813954277979023
$arr[16,10] = 0.09068194453499999
$arr[16,11] = 0.8161375008149999
$arr[16,12] = 0.0006625986182200453
$arr[16,13] = 0.0007868638160796272
$arr[17,0] = 0.010485
$arr[17,1] = 0.031455
$arr[17,2] = 0.001624843498075486
$arr[17,3] = 0.001634547755634201
$arr[17,4] = 2
$arr[17,5] = 0.6666666666666666
$arr[17,6] = 0.637617
$arr[17,7] = 1.912851
$arr[17,8] = 0.03006398091472189
$arr[17,9] = 0.03549028273468269
$arr[17,10] = 0.006685414244999999
$arr[17,11] = 0.06016872820499999
$arr[17,12] = 0.00004884926391555138
$arr[17,13] = 0.00005801056199079883
$arr[18,0] = 0.010485
$arr[18,1] = 0.031455
$arr[18,2] = 0.001624843498075486
$arr[18,3] = 0.001634547755634201
$arr[18,4] = 2
$arr[18,5] = 1
$arr[18,6] = 9.728125
$arr[18,7] = 19.45625
$arr[18,8] = 0.4586862714388558
$arr[18,9] = 0.3609835859963323
$arr[18,10] = 0.101999390625
$arr[18,11] = 0.61199634375
$arr[18,12] = 0.0007452934058039126
$arr[18,13] = 0.0005900449103110906
$arr[19,0] = 0.010485
$arr[19,1] = 0.031455
$arr[19,2] = 0.001624843498075486
$arr[19,3] = 0.001634547755634201
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 0.8878186666666666
$arr[19,7] = 2.663456
$arr[19,8] = 0.0418611226651744
$arr[19,9] = 0.0494167117519279
$arr[19,10] = 0.009308778719999999
$arr[19,11] = 0.08377900847999999
$arr[19,12] = 0.00006801777298464901
$arr[19,13] = 0.00008077397528493599
$ws.Range("G2:T21").Value = $arr
